$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.7810574977206343
$ws.Range("C2").Value = -0.6598279200072594
$ws.Range("D2").Value = -0.7874574027300708

$ws.Range("B3").Value = 0.6690545589678621
$ws.Range("C3").Value = -0.5021097860384722
$ws.Range("D3").Value = -0.5264775849324478

$ws.Range("B4").Value = -0.7561300374331079
$ws.Range("C4").Value = -0.6984879414044283
$ws.Range("D4").Value = -0.8016463904492847

$ws.Range("B5").Value = 0.7567730340215171
$ws.Range("C5").Value = 0.5452422622524453
$ws.Range("D5").Value = 0.5335943917643384

$ws.Range("B6").Value = 0.5434166752942166
$ws.Range("C6").Value = -0.5402719076678451
$ws.Range("D6").Value = 0.5362581232391506

$ws.Range("B7").Value = -0.8369907603175941
$ws.Range("C7").Value = -0.5939756606374655
$ws.Range("D7").Value = -0.687149874758077

$ws.Range("B8").Value = 0.6853874840347415
$ws.Range("C8").Value = -0.6437703782022719
$ws.Range("D8").Value = -0.5434394135265664

$ws.Range("B9").Value = 0.8304173421622492
$ws.Range("C9").Value = -0.7618437608136642
$ws.Range("D9").Value = 0.7666860027976375
